$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "Backend Java Developer"
$ws.Range("B80").Value = "https://www.dice.com/job-detail/480135e6-b37f-48ef-b352-ad7c9c052ba9"
$ws.Range("C80").Value = "Minneapolis, Minnesota"
$ws.Range("D80").Value = "Contract"
$ws.Range("E80").Value = "USD 50.00 - 60.00 per hour"
$ws.Range("F80").Value = "Robert Half"

$ws.Range("A81").Value = "Backend Java Developer"
$ws.Range("B81").Value = "https://www.dice.com/job-detail/b478024f-0940-4af2-8f3c-bb021830c8f7"
$ws.Range("C81").Value = "Minneapolis, Minnesota"
$ws.Range("D81").Value = "Contract"
$ws.Range("E81").Value = "USD 60.00 - 70.00 per hour"
$ws.Range("F81").Value = "PETADATA"
